$d = $word.ActiveDocument
$sec = $d.Sections.First

# Rename the Pearson Edexcel logo pictures (currently "image1.png") to
# "image2.png" in both the first-page footer and the default footer.
for ($i = 1; $i -le 3; $i++) {
    $ftr = $sec.Footers.Item($i)
    if ($ftr.Exists) {
        $shapes = $ftr.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $shp = $shapes.Item($j)
            if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                $shp.Name = "image2.png"
            }
        }
    }
}

# Rename the BTEC logo picture (currently "image2.jpg") to "image1.jpg"
# in the first-page header.
for ($i = 1; $i -le 3; $i++) {
    $hdr = $sec.Headers.Item($i)
    if ($hdr.Exists) {
        $shapes = $hdr.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $shp = $shapes.Item($j)
            if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                $shp.Name = "image1.jpg"
            }
        }
    }
}
